# Auto-committed edit: add a new "MrKey" (交易編號) field row to the
# TxAuthorize DB layout sheet, renumber the SEQ column accordingly, and
# leave the selection where the author's session ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# --- Insert a new row above the old row 17 (the "AutoSeq / KeyID" row),
#     pushing the remaining rows (old 17-20) down to 18-21. ---
$ws.Rows("17:17").Insert()

# Copy the formatting of the row above (row 16, the last "varchar2 field"
# row) onto the freshly inserted row 17 so it matches the rest of the table.
$ws.Range("A16:G16").Copy()
$ws.Range("A17:G17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Fill in the new field row: SEQ 9, MrKey / 交易編號, varchar2(20). ---
$ws.Range("A17").Value = 9
$ws.Range("B17").Value = "MrKey"
$ws.Range("C17").Value = "交易編號"
$ws.Range("D17").Value = "varchar2"
$ws.Range("E17").Value = "20"

# --- Renumber the SEQ column for the rows above the insertion point so the
#     sequence stays contiguous (5,6,7,8,9 -> 4,5,6,7,8) now that the new
#     row takes SEQ 9. ---
$ws.Range("A12").Value = 4
$ws.Range("A13").Value = 5
$ws.Range("A14").Value = 6
$ws.Range("A15").Value = 7
$ws.Range("A16").Value = 8

# --- Leave the selection/viewport where the author ended up editing. ---
$ws.Range("G23").Select()
